$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1. Solver parameters (defined names) - point at the new second problem block
# ---------------------------------------------------------------------------
$wb.Names.Item("solver_adj").RefersTo  = "=Arkusz1!`$A`$16:`$B`$16"
$wb.Names.Item("solver_lhs1").RefersTo = "=Arkusz1!`$K`$16"
$wb.Names.Item("solver_lhs2").RefersTo = "=Arkusz1!`$K`$17"
$wb.Names.Item("solver_lhs3").RefersTo = "=Arkusz1!`$K`$18"
$wb.Names.Item("solver_lhs4").RefersTo = "=Arkusz1!`$K`$19"
$wb.Names.Item("solver_lhs5").RefersTo = "=Arkusz1!`$K`$19"
$wb.Names.Item("solver_num").RefersTo  = "=3"
$wb.Names.Item("solver_opt").RefersTo  = "=Arkusz1!`$E`$16"
$wb.Names.Item("solver_rel1").RefersTo = "=3"
$wb.Names.Item("solver_rel2").RefersTo = "=3"
$wb.Names.Item("solver_rel3").RefersTo = "=3"
$wb.Names.Item("solver_rhs1").RefersTo = "=Arkusz1!`$J`$16"
$wb.Names.Item("solver_rhs2").RefersTo = "=Arkusz1!`$J`$17"
$wb.Names.Item("solver_rhs3").RefersTo = "=Arkusz1!`$J`$18"
$wb.Names.Item("solver_rhs4").RefersTo = "=Arkusz1!`$J`$19"
$wb.Names.Item("solver_rhs5").RefersTo = "=Arkusz1!`$J`$19"
$wb.Names.Item("solver_typ").RefersTo  = "=2"

# ---------------------------------------------------------------------------
# 2. New "b)" task block, rows 13-21 (mirrors the "a)" block in rows 1-10)
# ---------------------------------------------------------------------------

# -- row 13 : section label ---------------------------------------------
$ws.Range("A13").Value = "b)"

# -- row 14 : headers -----------------------------------------------------
$ws.Range("A14").Value = "zmienne decyzyjne:"
$ws.Range("D14").Value = "funkcja celu:"
$ws.Range("G14").Value = "zbiór dopuszczalny:"

# -- row 15 : sub headers (bordered) --------------------------------------
$ws.Range("A15").Value = "x"
$ws.Range("B15").Value = "y"
$ws.Range("D15").Value = "f(x,y)=ax+by"

$ws.Range("G15:K15").Borders.Color = 0
$ws.Range("G15:K15").Borders.LineStyle = 1
$ws.Range("G15").Value = "przy x "
$ws.Range("H15").Value = "przy y"
$ws.Range("I15").Value = "znak"
$ws.Range("J15").Value = "ograniczenie"
$ws.Range("K15").Value = "formuła"

# -- row 16 : a/b coefficients + objective + first constraint row --------
$ws.Range("A16:B16").Interior.Color = 65535
$ws.Range("A16:B16").Borders.Color = 0
$ws.Range("A16:B16").Borders.LineStyle = 1
$ws.Range("A16").Value = 2.0000000000000004
$ws.Range("B16").Value = 0.49999999999999978

$ws.Range("D16").Value = "f(x,y)="

$ws.Range("E16").Interior.Color = 65535
$ws.Range("E16").Formula = "=SUMPRODUCT(A16:B16,A18:B18)"

$ws.Range("G16:K16").Borders.Color = 0
$ws.Range("G16:K16").Borders.LineStyle = 1
$ws.Range("G16").Value = 1
$ws.Range("H16").Value = 2
$ws.Range("I16").Value = ">="
$ws.Range("J16").Value = 3
$ws.Range("K16").Formula = "=SUMPRODUCT(`$A`$16:`$B`$16,G16:H16)"

# -- row 17 : a,b labels + "min" + second constraint row ------------------
$ws.Range("A17").Value = "a"
$ws.Range("B17").Value = "b"

$ws.Range("D17").Font.Color = 255
$ws.Range("D17").Value = "min"

$ws.Range("G17:K17").Borders.Color = 0
$ws.Range("G17:K17").Borders.LineStyle = 1
$ws.Range("G17").Value = 1
$ws.Range("H17").Value = 4
$ws.Range("I17").Value = ">="
$ws.Range("J17").Value = 4
$ws.Range("K17").Formula = "=SUMPRODUCT(`$A`$16:`$B`$16,G17:H17)"

# -- row 18 : constraint coefficients (bordered) + third constraint row ---
$ws.Range("A18:B18").Borders.Color = 0
$ws.Range("A18:B18").Borders.LineStyle = 1
$ws.Range("A18").Value = 12
$ws.Range("B18").Value = 42

$ws.Range("G18:K18").Borders.Color = 0
$ws.Range("G18:K18").Borders.LineStyle = 1
$ws.Range("G18").Value = 3
$ws.Range("H18").Value = 1
$ws.Range("I18").Value = ">="
$ws.Range("J18").Value = 3
$ws.Range("K18").Formula = "=SUMPRODUCT(`$A`$16:`$B`$16,G18:H18)"

# -- row 19 : fourth constraint row ----------------------------------------
$ws.Range("G19:K19").Borders.Color = 0
$ws.Range("G19:K19").Borders.LineStyle = 1
$ws.Range("G19").Value = 1
$ws.Range("H19").Value = 0
$ws.Range("I19").Value = ">="
$ws.Range("J19").Value = 0
$ws.Range("K19").Formula = "=SUMPRODUCT(`$A`$16:`$B`$16,G19:H19)"

# -- row 20 : fifth constraint row -----------------------------------------
$ws.Range("G20:K20").Borders.Color = 0
$ws.Range("G20:K20").Borders.LineStyle = 1
$ws.Range("G20").Value = 0
$ws.Range("H20").Value = 1
$ws.Range("I20").Value = ">="
$ws.Range("J20").Value = 0
$ws.Range("K20").Formula = "=SUMPRODUCT(`$A`$16:`$B`$16,G20:H20)"

# -- row 21 : trailing bordered blank row ----------------------------------
$ws.Range("G21:K21").Borders.Color = 0
$ws.Range("G21:K21").Borders.LineStyle = 1

# ---------------------------------------------------------------------------
# 3. Selection moves to the new objective cell
# ---------------------------------------------------------------------------
[void]$ws.Range("E16").Select()
